# "add lock screen in step 1 for mitigate mutiples requests"
# The visible workbook change is the document date stamp in A1 advancing
# from 2024-02-22 (serial 45344) to 2024-04-23 (serial 45405).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")  # "CINTA EMBALAR" sheet (first/active tab)

# A1 holds the printed document date as a real date serial (numFmt 14).
# Assign an actual Date so it round-trips as a serial number (45405),
# not a formatted string.
$ws.Range("A1").Value = (Get-Date -Year 2024 -Month 4 -Day 23 -Hour 0 -Minute 0 -Second 0).Date
